# Jogos_da_Semana_FlashScore_FULL_2024-11-22.xlsx - odds refresh
# Updates a handful of odds cells in rows 3, 5, 7 and 17, then removes the two
# stale fixtures that used to live in rows 18-19 (row 20's fixture shifts up
# to become the new row 18, and the sheet shrinks from A1:BD20 to A1:BD18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 odds updates ---
$ws.Range("J3").Value = 2.05
$ws.Range("S3").Value = 1.33
$ws.Range("T3").Value = 3.25
$ws.Range("U3").Value = 1.8
$ws.Range("V3").Value = 1.91
$ws.Range("AD3").Value = 7.5
$ws.Range("AF3").Value = 51
$ws.Range("AS3").Value = 126
$ws.Range("AT3").Value = 3.25
$ws.Range("AU3").Value = 8.5
$ws.Range("BB3").Value = 126

# --- Row 5 odds updates ---
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 9
$ws.Range("Q5").Value = 2.1
$ws.Range("R5").Value = 1.7

# --- Row 7 odds updates ---
$ws.Range("Q7").Value = 2.03
$ws.Range("R7").Value = 1.83

# --- Row 17 odds updates ---
$ws.Range("G17").Value = 3.3
$ws.Range("H17").Value = 4.3
$ws.Range("I17").Value = 1.75
$ws.Range("J17").Value = 3.45
$ws.Range("K17").Value = 2.72
$ws.Range("L17").Value = 2.12
$ws.Range("W17").Value = 19.5
$ws.Range("X17").Value = 24
$ws.Range("Y17").Value = 11.5
$ws.Range("Z17").Value = 40
$ws.Range("AA17").Value = 19.5
$ws.Range("AB17").Value = 17.5
$ws.Range("AC17").Value = 28
$ws.Range("AD17").Value = 9.5
$ws.Range("AE17").Value = 10
$ws.Range("AF17").Value = 22
$ws.Range("AG17").Value = 75
$ws.Range("AH17").Value = 14
$ws.Range("AI17").Value = 12.5
$ws.Range("AK17").Value = 15.5
$ws.Range("AL17").Value = 10.25
$ws.Range("AM17").Value = 12
$ws.Range("AN17").Value = 6.5
$ws.Range("AO17").Value = 16
$ws.Range("AP17").Value = 15.5
$ws.Range("AQ17").Value = 60
$ws.Range("AR17").Value = 55
$ws.Range("AS17").Value = 100
$ws.Range("AT17").Value = 4.9
$ws.Range("AU17").Value = 6.2
$ws.Range("AX17").Value = 4.65
$ws.Range("AY17").Value = 8.25
$ws.Range("AZ17").Value = 11
$ws.Range("BA17").Value = 21
$ws.Range("BB17").Value = 28

# --- Remove the two fixtures that used to sit in rows 18 and 19 ---
# (Sporting CP - Amarante, and Boston River - Defensor Sp.). Deleting row 18
# twice pulls the old row 19 and then the old row 20 fixture up in turn, so
# the Nacional - Montevideo City match ends up as the new row 18 and the
# sheet's used range shrinks to A1:BD18.
$ws.Rows.Item(18).Delete()
$ws.Rows.Item(18).Delete()
